$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used inside "Virat Kohli (c)" in the source data.
$nbsp = [char]0x00A0
$batsman = "Virat Kohli" + $nbsp + "(c)"

# Insert two new columns (ownTeam, oppTeam) before the existing "batsman" column (D),
# shifting batsman..sr from D:I to F:K.
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "venue"
$ws.Cells.Item(1,2).Value = "date"
$ws.Cells.Item(1,3).Value = "result"
$ws.Cells.Item(1,4).Value = "ownTeam"
$ws.Cells.Item(1,5).Value = "oppTeam"
$ws.Cells.Item(1,6).Value = "batsman"
$ws.Cells.Item(1,7).Value = "totalRuns"
$ws.Cells.Item(1,8).Value = "totalBalls"
$ws.Cells.Item(1,9).Value = "total4s"
$ws.Cells.Item(1,10).Value = "total6s"
$ws.Cells.Item(1,11).Value = "sr"

# ---- Data rows (venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr) ----
$rows = @(
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", $batsman, "9", "14", "0", "0", "64.28"),
    @(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Rajasthan Royals", $batsman, "43", "32", "1", "2", "134.37"),
    @(" Abu Dhabi", " October 03 2020", "RCB won by 8 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Rajasthan Royals", $batsman, "72", "53", "7", "2", "135.84"),
    @(" Dubai (DSC)", " September 24 2020", "Kings XI won by 97 runs", "Royal Challengers Bangalore", "Kings XI Punjab", $batsman, "1", "5", "0", "0", "20.00"),
    @(" Sharjah", " October 15 2020", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", $batsman, "48", "39", "3", "0", "123.07"),
    @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "14", "13", "0", "0", "107.69"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "6", "7", "0", "0", "85.71"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", $batsman, "43", "39", "2", "1", "110.25"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "7", "7", "0", "0", "100.00"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", $batsman, "29", "24", "2", "1", "120.83"),
    @(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Royal Challengers Bangalore", "Kolkata Knight Riders", $batsman, "18", "17", "2", "0", "105.88"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Royal Challengers Bangalore", "Mumbai Indians", $batsman, "3", "11", "0", "0", "27.27"),
    @(" Sharjah", " October 12 2020", "RCB won by 82 runs", "Royal Challengers Bangalore", "Kolkata Knight Riders", $batsman, "33", "28", "1", "0", "117.85"),
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", $batsman, "50", "43", "1", "1", "116.27"),
    @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", $batsman, "90", "52", "4", "4", "173.07")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $ws.Cells.Item($r,7).Value = "'" + $row[6]
    $ws.Cells.Item($r,8).Value = "'" + $row[7]
    $ws.Cells.Item($r,9).Value = "'" + $row[8]
    $ws.Cells.Item($r,10).Value = "'" + $row[9]
    $ws.Cells.Item($r,11).Value = "'" + $row[10]
    $r = $r + 1
}
